$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11 (shifts existing rows 11..104 down to 12..105)
$ws.Rows("11:11").Insert()

# Populate the new row 11 with the FxE "output" / "configuration_fxe" entry
$ws.Range("A11").Value = "CHE"
$ws.Range("B11").Value = "trd_gas"
$ws.Range("C11").Value = "output"
$ws.Range("C11").VerticalAlignment = -4108
$ws.Range("D11").Value = "configuration_fxe"
$ws.Range("F11").Value = "gas"
$ws.Range("G11").Value = 1

# Re-apply AutoFilter over the now-larger data range (A5:L303 -> A5:L304)
$ws.AutoFilterMode = $false
$ws.Range("A5:L304").AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new range
$fdb = $wb.Names.Item("Sheet1!_FilterDatabase")
$fdb.RefersTo = "=Sheet1!`$A`$5:`$L`$304"

# Update selection to match the saved view state
$ws.Range("G12").Select()
